# 2021-07 Victoria Outbreak Paths.xlsx - "Add files via upload" commit
# Re-creates the edit described by the OOXML diff:
#  - On "Sheet1": a handful of existing Sub-Cluster / Link Label cells get
#    corrected text ("MCG" -> "AAMI Park", and one "Ms Frankie Restaurant
#    Patron" -> "...Close Contacts"), and 12 brand-new rows (150-161, all
#    dated 24-Jul-2021 / serial 44401) are appended to Table1.
#  - On "Date Colours": the gradient of day colours in column B is
#    regenerated (each existing row gets a new shade) and a 13th row/shade
#    is added for the new day.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Date Colours")

# ---------------------------------------------------------------------
# 1. Fix existing Sub-Cluster/Link Label text on Sheet1
# ---------------------------------------------------------------------
foreach ($addr in @("E54","E55","E70","E102","E103","E135","E137","E138")) {
    $ws1.Range($addr).Value = "AAMI Park"
}
$ws1.Range("F141").Value = "Ms Frankie Restaurant Patron Close Contacts"

# ---------------------------------------------------------------------
# 2. Append the 12 new outbreak-path rows to Sheet1 (rows 150-161)
# ---------------------------------------------------------------------
$newRows = @(
    @(150, "M44",     "M133",                       "Maribyrnong", "AAMI Park",                  "AAMI Park Close Contacts",                   "Delta (B.1.617.2)", "Isolated"),
    @(151, "M119 m",  "M134",                       "Maribyrnong", "AAMI Park",                  "Household",                                   "Delta (B.1.617.2)", "Isolated"),
    @(152, "M89",     "M135",                       "Maribyrnong", "AAMI Park",                  "Household",                                   "Delta (B.1.617.2)", "Isolated"),
    @(153, "M120",    "M136",                       "Maribyrnong", "AAMI Park",                  "Household",                                   "Delta (B.1.617.2)", "Isolated"),
    @(154, "M29",     "M137",                       "Maribyrnong", "AAMI Park",                  "Household",                                   "Delta (B.1.617.2)", "Isolated"),
    @(155, "M63",     "M138",                       "Maribyrnong", "Ms Frankie Restaurant",      "Ms Frankie Restaurant Patron Close Contacts", "Delta (B.1.617.2)", "Isolated"),
    @(156, "M80",     "M139",                       "Maribyrnong", "Ms Frankie Restaurant",      "Household",                                   "Delta (B.1.617.2)", "Isolated"),
    @(157, "M80",     "M140",                       "Maribyrnong", "Ms Frankie Restaurant",      "Household",                                   "Delta (B.1.617.2)", "Isolated"),
    @(158, "M141",    "M141",                       "Maribyrnong", "Isola Apartments, Hawthorn", "Isola Apartments, Hawthorn",                  "Delta (B.1.617.2)", "Wild"),
    @(159, "M141",    "M142",                       "Maribyrnong", "Isola Apartments, Hawthorn", "Isola Apartments, Hawthorn",                  "Delta (B.1.617.2)", "Wild"),
    @(160, "M20",     "M143 10s",                   "Maribyrnong", "Bacchus Marsh Grammar",      "Bacchus Marsh Grammar Student",               "Delta (B.1.617.2)", "Isolated"),
    @(161, "M86",     "M144",                       "Maribyrnong", "Young & Jacksons",           "Household",                                   "Delta (B.1.617.2)", "Isolated")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws1.Range("A$r").Value = 44401
    $ws1.Range("A$r").NumberFormat = $ws1.Range("A149").NumberFormat
    $ws1.Range("B$r").Value = $row[1]
    $ws1.Range("C$r").Value = $row[2]
    $ws1.Range("D$r").Value = $row[3]
    $ws1.Range("E$r").Value = $row[4]
    $ws1.Range("F$r").Value = $row[5]
    $ws1.Range("G$r").Value = $row[6]
    $ws1.Range("H$r").Value = $row[7]
}

# Grow Table1 so it covers the new rows, keeping the header row intact
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:H161"))

# Move the selection the same way Excel leaves it after typing into H161
$ws1.Range("H161").Select()

# ---------------------------------------------------------------------
# 3. Regenerate the "Colour Code" gradient on the Date Colours sheet and
#    add the extra day (row 13)
# ---------------------------------------------------------------------
$colourCodes = @(
    "#f9f8fc",
    "#f2f2f9",
    "#ecebf5",
    "#e5e5f2",
    "#dfdeef",
    "#d8d8ec",
    "#d2d1e8",
    "#cccbe5",
    "#c5c5e2",
    "#bfbedf",
    "#b8b8db",
    "#b2b2d8"
)

for ($i = 0; $i -lt $colourCodes.Length; $i++) {
    $r = $i + 2
    $ws2.Range("B$r").Value = $colourCodes[$i]
}
